# Generate Report for Handback
#
# The 77a4a6e8-9225-47f5-b53f-bfe9c910a317 source file has been handed
# back from localization ("in sync with en-US" - no content changes were
# required), so its status flips from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it's reported, and the
# per-language detail sheets gain the handback bookkeeping columns
# (Latest Target File / Latest Handback File / Latest Handback DateTime).

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: update the status columns for the 77a4a6e8 file row.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusHandedBack
$overview.Range("C2").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusHandedBack
$zhcn.Range("H2").Value = "2016-03-22 00:34:21"

$zhcn.Hyperlinks.Add(
    $zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/bfbb6bf5e628a710fedd92ee5e3fe32034f8d1e9/e2e/77a4a6e8-9225-47f5-b53f-bfe9c910a317.md",
    "",
    "",
    "77a4a6e8-9225-47f5-b53f-bfe9c910a317.md"
) | Out-Null

$zhcn.Hyperlinks.Add(
    $zhcn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0c7b4e99afdc51ba53d615aa2fb31226115ae1f9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/77a4a6e8-9225-47f5-b53f-bfe9c910a317.61be0d3a06f1e2c0b0a43f3cebec17b440913912.zh-cn.xlf",
    "",
    "",
    "77a4a6e8-9225-47f5-b53f-bfe9c910a317.61be0d3a06f1e2c0b0a43f3cebec17b440913912.zh-cn.xlf"
) | Out-Null

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusHandedBack
$dede.Range("H2").Value = "2016-03-22 00:34:28"

$dede.Hyperlinks.Add(
    $dede.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/bfbb6bf5e628a710fedd92ee5e3fe32034f8d1e9/e2e/77a4a6e8-9225-47f5-b53f-bfe9c910a317.md",
    "",
    "",
    "77a4a6e8-9225-47f5-b53f-bfe9c910a317.md"
) | Out-Null

$dede.Hyperlinks.Add(
    $dede.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/833c5f2d2b94ddbac60635c343076cc223d768ed/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/77a4a6e8-9225-47f5-b53f-bfe9c910a317.61be0d3a06f1e2c0b0a43f3cebec17b440913912.de-de.xlf",
    "",
    "",
    "77a4a6e8-9225-47f5-b53f-bfe9c910a317.61be0d3a06f1e2c0b0a43f3cebec17b440913912.de-de.xlf"
) | Out-Null
